$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A6").Value = "Demo inplannen"
$ws.Range("B6").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("D6").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("F6").Value = "2025-08-19 19:39:13"
$ws.Range("G6").Value = "Nee"
$ws.Range("H6").Value = "Ja"
$ws.Range("I6").Value = "Nee"
$ws.Range("J6").Value = "Nee"

$ws.Range("D2:D5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D6"))
$ws.Range("G2:G5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G6"))
$ws.Range("H2:H5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H6"))
$ws.Range("I2:I5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I6"))
$ws.Range("J2:J5").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J6"))

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 5
